# Auto-generated edit script: updates crypto price/volume table cells
# per the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.299.32"
$ws.Range("E2").Value = "  +2.32%  "
$ws.Range("D3").Value = "3.579.30"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Formula = "'241.17"
$ws.Range("E5").Value = "  +2.62%  "
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").Formula = "'1.73"
$ws.Range("E6").Value = "  +17.16%  "
$ws.Range("B7").Value = "BNB"
$ws.Range("C7").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D7").Formula = "'655.84"
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("D8").Formula = "'0.426"
$ws.Range("E8").Value = "  +6.88%  "
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("E10").Value = "  +4.89%  "
$ws.Range("D11").Value = "3.577.63"
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").Formula = "'44.22"
$ws.Range("E12").Value = "  +4.37%  "
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").Value = "4.242.98"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").Value = "96.853.90"
$ws.Range("E16").Value = "  +1.93%  "
$ws.Range("E17").Value = "  +3.05%  "
$ws.Range("D18").Formula = "'8.66"
$ws.Range("E18").Value = "  +11.66%  "
$ws.Range("D19").Value = "3.573.28"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Formula = "'12.72"
$ws.Range("D21").Formula = "'18.04"
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("E22").Value = "  +10.25%  "
$ws.Range("D23").Formula = "'3.49"
$ws.Range("E23").Value = "  +0.87%  "
$ws.Range("D24").Formula = "'514.19"
$ws.Range("E24").Value = "  +1.10%  "
$ws.Range("D25").Formula = "'0.0000206"
$ws.Range("E25").Value = "  +5.60%  "
$ws.Range("D26").Formula = "'6.87"
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("D27").Formula = "'101.80"
$ws.Range("E27").Value = "  +6.95%  "
$ws.Range("D28").Formula = "'13.07"
$ws.Range("E28").Value = "  +2.96%  "
$ws.Range("D29").Value = "3.769.60"
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("E30").Value = "  +17.57%  "
$ws.Range("E31").Value = "  -1.64%  "
$ws.Range("D32").Formula = "'11.93"
$ws.Range("E32").Value = "  +3.83%  "
$ws.Range("D33").Formula = "'0.999"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").Value = "  +4.02%  "
$ws.Range("D35").Formula = "'1.00"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Formula = "'31.84"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").Formula = "'8.85"
$ws.Range("E37").Value = "  +4.03%  "
$ws.Range("D38").Formula = "'615.71"
$ws.Range("E38").Value = "  +5.51%  "
$ws.Range("D39").Formula = "'0.565"
$ws.Range("E39").Value = "  +1.36%  "
$ws.Range("D40").Formula = "'1.64"
$ws.Range("E40").Value = "  -2.45%  "
$ws.Range("D41").Formula = "'1.95"
$ws.Range("E41").Value = "  +7.27%  "
$ws.Range("D42").Formula = "'0.155"
$ws.Range("E42").Value = "  +2.54%  "
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Formula = "'0.923"
$ws.Range("E44").Value = "  +1.91%  "
$ws.Range("E45").Value = "  +5.08%  "
$ws.Range("D46").Formula = "'0.0439"
$ws.Range("E46").Value = "  +6.05%  "
$ws.Range("D47").Formula = "'2.31"
$ws.Range("E47").Value = "  +1.73%  "
$ws.Range("D48").Formula = "'23.61"
$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("E49").Value = "  +32.28%  "
$ws.Range("E50").Value = "  +4.20%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Formula = "'32.98"
$ws.Range("E51").Value = "  -3.83%  "
